# Edit script for "variables and scoping.docx"
$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark from the end of the first paragraph
if ($d.Bookmarks.Exists("_GoBack")) {
    $b = $d.Bookmarks.Item("_GoBack")
    $b.Delete()
}

# 2. Merge the runs (and drop proofErr markers) in the "Because behavior differs..." paragraph
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
$r3trim = $d.Range($r3.Start, $r3.End - 1)
$r3trim.Text = "PLACEHOLDER_TEXT_MERGE_1"
$r3new = $d.Range($r3.Start, $r3.Start + 24)
$r3new.Text = "Because behavior differs between var and let, you can’t blindly replace var with let or const in existing code; you have to be careful during refactoring"

# 3. Merge the runs (and drop proofErr markers) in the "In ECMAScript 6..." paragraph
$p14 = $d.Paragraphs.Item(14)
$r14 = $p14.Range
$r14trim = $d.Range($r14.Start, $r14.End - 1)
$r14trim.Text = "PLACEHOLDER_TEXT_MERGE_2"
$r14new = $d.Range($r14.Start, $r14.Start + 24)
$r14new.Text = "In ECMAScript 6, accessing a let or const variable before its declaration (within its scope) causes a ReferenceError . The time span when that happens, between the creation of a variable's binding and its declaration, is called the temporal dead zone."

# 4. Append the new paragraphs (and re-add the _GoBack bookmark on the last one)
$pLast = $d.Paragraphs.Last
$endPos = $pLast.Range.End
$target = $d.Range($endPos, $endPos)
$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Lets allow for mutation</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Consts are immutable to an extent</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:tab/>
        <w:t>You cannot overwrite the very base of the const</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:tab/>
        <w:t>You can’t do something like const someConst = {}; someConst=’abc’;</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>That will throw a TypeError because of the strict binding that is produced for const-declared variables</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>However, you can do const someConst = {}; someConst.prop=123</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>Const only means that a variable always has the same value – it doesn’t mean that the value itself is or becomes immutable</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:ind w:left="720"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Temporal Dead Zones</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>When entering the scope of a let or a const, it can’t be accessed (got or set) until execution reaches the declaration</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>For vars, as soon as the scope for a var variable is entered, storage space is created for it, and the variable is immediately initialized, by setting it to undefined</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Only when the execution within the scope reaches the declaration, the variable is set to the value specified by the initializer </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>For lets, when the scope of a let (its surrounding block) is entered, storage space is created for it</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>However, it remains uninitialized – not even undefined is set as its value</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>this difference generates the temporal dead zone</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>when the execution within the scope reaches the declaration, the variable is set to the value specified by the initializer – if there is one – if there isn’t then the value of the variable is set to undefined</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>so the difference comes down to when the variable is initialized and the scoping (block versus function)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Why is there a temporal dead zone</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>There are several reasons why const and let have temporal dead zones</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Catch programming errors  - being able to access a variable before its declaration is strange</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Making const work properly is difficult – TDZs provide a rational semantics for const. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Future proofing for guards – a mechanism for enforcing at runtime that a variable has the correct value, hoisting generating a value of undefined may be in conflict with the guarantee given by its guard</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xmlFrag)
